$wb = $excel.ActiveWorkbook

# Switch to the Users sheet, update the name, and select A2
$ws = $wb.Worksheets.Item("Users")
$ws.Activate()
$ws.Range("A2").Value = "Thomas Bailey"
$ws.Range("A2").Select()
